$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 581-582 (shifts the old 581.. down to 583..)
$ws.Rows("581:582").Insert()

# Row 581 - new weekly record (copy shape of the row that used to sit here,
# with the updated market figures from this week's publication)
$ws.Range("A581").Value = 9
$ws.Range("B581").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C581").Value = "Metropolitana"
$ws.Range("D581").Value = 44714
$ws.Range("E581").Value = 13
$ws.Range("F581").Value = 100112040
$ws.Range("G581").Value = "Cilantro"
$ws.Range("H581").Value = "Sin especificar"
$ws.Range("I581").Value = "Primera"
$ws.Range("J581").Value = 52
$ws.Range("K581").Value = 5000
$ws.Range("L581").Value = 5000
$ws.Range("M581").Value = 5000
$ws.Range("N581").Value = "$/caja 36 atados"
$ws.Range("O581").Value = "Región Metropolitana"
$ws.Range("P581").Value = 139
$ws.Range("Q581").Value = 36
$ws.Range("R581").Value = "Hortaliza"

# Row 582 - new weekly record
$ws.Range("A582").Value = 9
$ws.Range("B582").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C582").Value = "Metropolitana"
$ws.Range("D582").Value = 44714
$ws.Range("E582").Value = 13
$ws.Range("F582").Value = 100112040
$ws.Range("G582").Value = "Cilantro"
$ws.Range("H582").Value = "Sin especificar"
$ws.Range("I582").Value = "Primera"
$ws.Range("J582").Value = 160
$ws.Range("K582").Value = 8000
$ws.Range("L582").Value = 9000
$ws.Range("M582").Value = 8500
$ws.Range("N582").Value = "$/docena de atados"
$ws.Range("O582").Value = "Región Metropolitana"
$ws.Range("P582").Value = 2833
$ws.Range("Q582").Value = 3
$ws.Range("R582").Value = "Hortaliza"

Write-Output "rows inserted and populated"
